$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 5 (license_id / field_license_wbddh / TRUE),
# shifting the remaining rows (format row) up.
$ws.Rows.Item(5).Delete()

# Select the new row 5 (previously row 6, now the last data row)
$ws.Range("A5:XFD5").Select()
